$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-09 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-10 Monday", 2)

$d.Content.Find.Execute("36×64=2304", $true, $false, $false, $false, $false, $true, 1, $false, "80×76=6080", 2)
$d.Content.Find.Execute("26×52=1352", $true, $false, $false, $false, $false, $true, 1, $false, "34×15=510", 2)
$d.Content.Find.Execute("48×44=2112", $true, $false, $false, $false, $false, $true, 1, $false, "82×71=5822", 2)
$d.Content.Find.Execute("95×49=4655", $true, $false, $false, $false, $false, $true, 1, $false, "55×65=3575", 2)
$d.Content.Find.Execute("94×51=4794", $true, $false, $false, $false, $false, $true, 1, $false, "76×12=912", 2)

$d.Content.Find.Execute("27×12=324", $true, $false, $false, $false, $false, $true, 1, $false, "90×35=3150", 2)
$d.Content.Find.Execute("65×22=1430", $true, $false, $false, $false, $false, $true, 1, $false, "26×25=650", 2)
$d.Content.Find.Execute("99×89=8811", $true, $false, $false, $false, $false, $true, 1, $false, "99×52=5148", 2)
$d.Content.Find.Execute("56×21=1176", $true, $false, $false, $false, $false, $true, 1, $false, "33×79=2607", 2)
$d.Content.Find.Execute("54×25=1350", $true, $false, $false, $false, $false, $true, 1, $false, "62×57=3534", 2)

$d.Content.Find.Execute("27×84=2268", $true, $false, $false, $false, $false, $true, 1, $false, "72×86=6192", 2)
$d.Content.Find.Execute("59×63=3717", $true, $false, $false, $false, $false, $true, 1, $false, "70×13=910", 2)
$d.Content.Find.Execute("33×73=2409", $true, $false, $false, $false, $false, $true, 1, $false, "48×91=4368", 2)
$d.Content.Find.Execute("14×45=630", $true, $false, $false, $false, $false, $true, 1, $false, "12×82=984", 2)
$d.Content.Find.Execute("82×67=5494", $true, $false, $false, $false, $false, $true, 1, $false, "49×23=1127", 2)

$d.Content.Find.Execute("18×93=1674", $true, $false, $false, $false, $false, $true, 1, $false, "74×54=3996", 2)
$d.Content.Find.Execute("15×69=1035", $true, $false, $false, $false, $false, $true, 1, $false, "49×11=539", 2)
$d.Content.Find.Execute("68×74=5032", $true, $false, $false, $false, $false, $true, 1, $false, "75×73=5475", 2)
$d.Content.Find.Execute("98×25=2450", $true, $false, $false, $false, $false, $true, 1, $false, "24×93=2232", 2)
$d.Content.Find.Execute("46×19=874", $true, $false, $false, $false, $false, $true, 1, $false, "99×94=9306", 2)

$d.Content.Find.Execute("59×79=4661", $true, $false, $false, $false, $false, $true, 1, $false, "80×43=3440", 2)
$d.Content.Find.Execute("25×60=1500", $true, $false, $false, $false, $false, $true, 1, $false, "83×80=6640", 2)
$d.Content.Find.Execute("46×12=552", $true, $false, $false, $false, $false, $true, 1, $false, "85×37=3145", 2)
$d.Content.Find.Execute("55×44=2420", $true, $false, $false, $false, $false, $true, 1, $false, "75×87=6525", 2)
$d.Content.Find.Execute("95×42=3990", $true, $false, $false, $false, $false, $true, 1, $false, "71×54=3834", 2)
